# ContosoLearn Competitor SWOT -> DAFO de competidores ContosoLearn
# English -> Spanish (es-ES) localization pass:
#  - translated title / SWOT headers ("Strengths:" etc.) / body copy
#  - body & header runs get explicit lang (es-ES / en-US / ar-SA) and
#    Aptos / Times New Roman fonts
#  - every paragraph gets an explicit left-to-right reading order
#    (serializes as <w:bidi w:val="0"/>)

$d = $word.ActiveDocument

function Set-EsEsLang($range) {
    $range.LanguageID = "es-ES"
    $range.LanguageIDFarEast = "en-US"
    $range.LanguageIDOther = "ar-SA"
}

function Set-AptosFont($range) {
    $range.Font.Name = "Aptos"
    $range.Font.NameFarEast = "Aptos"
    $range.Font.NameBi = "Times New Roman"
}

function Format-BodyRange($range) {
    $range.Bold = 0
    $range.BoldBi = 0
    $range.Italic = 0
    $range.ItalicBi = 0
    Set-AptosFont $range
    Set-EsEsLang $range
}

function Format-HeaderRange($range) {
    $range.Bold = 1
    $range.BoldBi = 1
    Set-AptosFont $range
    Set-EsEsLang $range
}

function Set-ParaLtr($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Format.ReadingOrder = 0
}

# Replace $oldText with $newText inside paragraph $paraIndex, then
# apply header/body formatting to exactly the replaced text (nothing
# else in the paragraph is touched). $kind is "header" or "body".
function Replace-InPara($paraIndex, $oldText, $newText, $kind) {
    $p = $d.Paragraphs.Item($paraIndex)
    $searchRange = $p.Range
    $find = $searchRange.Find
    $find.ClearFormatting()
    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if ($kind -eq "header") {
        Format-HeaderRange $searchRange
    } else {
        Format-BodyRange $searchRange
    }
}

# --- Paragraph 1: title ---
$old = "ContosoLearn Competitor SWOT"
$new = "DAFO de competidores ContosoLearn"
Replace-InPara 1 $old $new "body"
Set-ParaLtr 1

# --- Paragraph 2: "Fabrikam Learning:" (text unchanged, formatting refreshed) ---
$p2 = $d.Paragraphs.Item(2)
Format-HeaderRange $p2.Range
Set-ParaLtr 2

# --- Paragraph 3: Fabrikam Strengths ---
Replace-InPara 3 "Strengths:" "Fortalezas:" "header"
$old = " Fabrikam Learning provides a comprehensive set of analytics and reporting tools. It ensures the continuous monitoring of teaching and learning activities, as well as pinpointing problematic areas that need to be addressed."
$new = " Fabrikam Learning proporciona un conjunto completo de herramientas de análisis e informes. Garantiza la supervisión continua de las actividades de enseñanza y aprendizaje, así como la identificación de áreas problemáticas que deben abordarse."
Replace-InPara 3 $old $new "body"
Set-ParaLtr 3

# --- Paragraph 4: Fabrikam Weaknesses ---
Replace-InPara 4 "Weaknesses:" "Puntos débiles:" "header"
$old = " While Fabrikam Learning has robust reporting capabilities, it might be overwhelming for some users due to its comprehensive nature."
$new = " aunque Fabrikam Learning tiene funcionalidades de informes sólidas, puede ser abrumador para algunos usuarios debido a su naturaleza completa."
Replace-InPara 4 $old $new "body"
Set-ParaLtr 4

# --- Paragraph 5: Fabrikam Opportunities ---
Replace-InPara 5 "Opportunities:" "Oportunidades:" "header"
$old = " There is a growing demand for personalized learning experiences and data-driven recommendations. Fabrikam Learning can leverage its robust analytics and reporting tools to meet this demand."
$new = " hay una creciente demanda de experiencias de aprendizaje personalizadas y recomendaciones controladas por datos. Fabrikam Learning puede aprovechar sus sólidas herramientas de análisis e informes para satisfacer esta demanda."
Replace-InPara 5 $old $new "body"
Set-ParaLtr 5

# --- Paragraph 6: Fabrikam Threats ---
Replace-InPara 6 "Threats:" "Amenazas:" "header"
$old = " The eLearning market is highly competitive with many players offering similar features. Fabrikam Learning needs to continuously innovate to stay ahead."
$new = " el mercado de eLearning es altamente competitivo con muchos jugadores que ofrecen características similares. Fabrikam Learning debe innovar continuamente para mantenerse a la vanguardia."
Replace-InPara 6 $old $new "body"
Set-ParaLtr 6

# --- Paragraph 7: "AdatumLearn:" (text unchanged, formatting refreshed) ---
$p7 = $d.Paragraphs.Item(7)
Format-HeaderRange $p7.Range
Set-ParaLtr 7

# --- Paragraph 8: AdatumLearn Strengths ---
Replace-InPara 8 "Strengths:" "Fortalezas:" "header"
$old = " AdatumLearn offers courses on business analysis techniques such as MOST and SWOT. This shows their commitment to providing valuable content to their users."
$new = " AdatumLearn ofrece cursos sobre técnicas de análisis de negocios como MOST y DAFO. Esto muestra su compromiso de proporcionar contenido valioso a sus usuarios."
Replace-InPara 8 $old $new "body"
Set-ParaLtr 8

# --- Paragraph 9: AdatumLearn Weaknesses ---
Replace-InPara 9 "Weaknesses:" "Puntos débiles:" "header"
$old = " The information provided in their courses is a compilation of third-party generated information. This might not be as valuable as original content."
$new = " la información proporcionada en sus cursos es una compilación de información generada por terceros. Esto podría no ser tan valioso como el contenido original."
Replace-InPara 9 $old $new "body"
Set-ParaLtr 9

# --- Paragraph 10: AdatumLearn Opportunities ---
Replace-InPara 10 "Opportunities:" "Oportunidades:" "header"
$old = " AdatumLearn can create more original content to provide unique value to their users. They can also expand their course offerings to cover more topics."
$new = " AdatumLearn puede crear contenido más original para proporcionar un valor único a sus usuarios. También puede ampliar sus ofertas de cursos para tratar más temas."
Replace-InPara 10 $old $new "body"
Set-ParaLtr 10

# --- Paragraph 11: AdatumLearn Threats ---
Replace-InPara 11 "Threats:" "Amenazas:" "header"
$old = " Like Fabrikam Learning, AdatumLearn also faces stiff competition in the eLearning market. They need to continuously improve their offerings to stay competitive.`""
$new = " al igual que Fabrikam Learning, AdatumLearn también se enfrenta a una competencia rígida en el mercado de eLearning. Necesita mejorar continuamente su oferta para mantenerse competitivo`"."
Replace-InPara 11 $old $new "body"
Set-ParaLtr 11

Write-Output "Done."
